# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff"
# - Timestamps bumped to reflect the new handoff generation time
# - Status column widened (content grew longer) on all three sheets

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-23 16:42:57"
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-23 16:42:52"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-23 16:42:57"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
